# "新题" sheet: log today's (2019-03-21) entry -- "55 dp", done.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("新题")

$ws.Range("A13").Value = [DateTime]"2019-03-21"
$ws.Range("A13").NumberFormat = "m/d/yy"   # reuse the existing short-date style (s=1), same as A2..A12
$ws.Range("B13").Value = "55 dp"
$ws.Range("E13").Value = "done"

$ws.Range("E13").Select()
